# "created excel management system"
#
# Sheet "user": re-label the columns (drop the leading "No." column,
# shifting every header one slot to the left) and add the first data row.
#
# Sheet "log": same left-shift of the header row (drop "No."), add a
# "Golongan" column, replace the old single sample/formula row with the
# real trip log rows, size a handful of columns, and move the selection.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "user"
# ---------------------------------------------------------------------
$wsUser = $wb.Worksheets.Item("user")

$wsUser.Range("A1").Value = "Mail"
$wsUser.Range("B1").Value = "Name"
$wsUser.Range("C1").Value = "ID"
$wsUser.Range("D1").Value = "Pass"
$wsUser.Range("E1").Value = "Saldo"
$wsUser.Range("F1").ClearContents()

$wsUser.Range("B2").Value = "bagas@mail.com"
$wsUser.Range("C2").Value = "Bagas"
# Looks numeric but must stay text - lead with an apostrophe like Excel does
# when a user types a quoted numeric string into a cell, then drop the
# resulting quote-prefix style so the cell is plain (unstyled) text.
$wsUser.Range("E2").Value = "'123456"
$wsUser.Range("E2").ClearFormats()
$wsUser.Range("F2").Value = 10000

# ---------------------------------------------------------------------
# Sheet "log"
# ---------------------------------------------------------------------
$wsLog = $wb.Worksheets.Item("log")

# Wipe the old placeholder row (incl. the ABS() formula in H2) so nothing
# stale survives the reshuffle below.
$wsLog.UsedRange.Clear()

$wsLog.Range("A1").Value = "Time In"
$wsLog.Range("B1").Value = "Gate In"
$wsLog.Range("C1").Value = "KM In"
$wsLog.Range("D1").Value = "Time Out"
$wsLog.Range("E1").Value = "Gate Out"
$wsLog.Range("F1").Value = "KM Out"
$wsLog.Range("G1").Value = "Jarak"
$wsLog.Range("H1").Value = "User"
$wsLog.Range("I1").Value = "Golongan"

$wsLog.Range("A2").Value = "11:35:33"
$wsLog.Range("B2").Value = "Semarang"
$wsLog.Range("D2").Value = "11:35:33"
$wsLog.Range("E2").Value = "Serpong"

$wsLog.Range("A3").Value = "11:38:45"
$wsLog.Range("B3").Value = "Juanda"

$wsLog.Range("D4").Value = "11:38:45"
$wsLog.Range("E4").Value = "Solo"
$wsLog.Range("H4").Value = "bagas@mail.com"

$wsLog.Range("A5").Value = "12:09:47"
$wsLog.Range("B5").Value = "Semarang"
$wsLog.Range("D5").Value = "12:09:47"
$wsLog.Range("E5").Value = "Taman Mini"
$wsLog.Range("H5").Value = "test2@mail.com"

$wsLog.Range("A6").Value = "12:18:25"
$wsLog.Range("B6").Value = "Tambak Sumur"
$wsLog.Range("D6").Value = "12:18:25"
$wsLog.Range("E6").Value = "Tanjung Priok"

$wsLog.Range("A7").Value = "12:23:34"
$wsLog.Range("B7").Value = "Tambak Oso"
$wsLog.Range("D7").Value = "12:23:34"
$wsLog.Range("E7").Value = "Juanda"
$wsLog.Range("H7").Value = "test1@mail.com"

$wsLog.Range("A8").Value = "12:37:13"
$wsLog.Range("B8").Value = "Bawen"
$wsLog.Range("C8").Value = 23.1
$wsLog.Range("D8").Value = "12:37:13"
$wsLog.Range("E8").Value = "Taman Mini"
$wsLog.Range("F8").Value = 4.5
$wsLog.Range("H8").Value = 473.6

$wsLog.Range("A9").Value = "12:44:40"
$wsLog.Range("B9").Value = "Bawen"
$wsLog.Range("D9").Value = "12:44:40"
$wsLog.Range("E9").Value = "Taman Mini"
$wsLog.Range("H9").Value = "test2@mail.com"

# Column widths (engine quantizes ColumnWidth to 1/6-char steps internally,
# so feed it the input that lands closest to the recorded OOXML width).
$wsLog.Columns.Item(1).ColumnWidth = 18.166666666666668
$wsLog.Columns.Item(2).ColumnWidth = 19.833333333333336
$wsLog.Columns.Item(4).ColumnWidth = 19.5
$wsLog.Columns.Item(5).ColumnWidth = 17.333333333333336
$wsLog.Columns.Item(8).ColumnWidth = 13.666666666666666

$wsLog.Activate() | Out-Null
$wsLog.Range("I3").Select() | Out-Null
